$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '67.604.65'
$ws.Range("E2").Value = '  +6.98%  '

Set-TextValue $ws.Range("D3") '3.505.83'
$ws.Range("E3").Value = '  +8.14%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.06%  '

Set-TextValue $ws.Range("D5") '191.01'
$ws.Range("E5").Value = '  +13.69%  '

Set-TextValue $ws.Range("D6") '554.57'
$ws.Range("E6").Value = '  +8.20%  '

Set-TextValue $ws.Range("D7") '3.499.06'
$ws.Range("E7").Value = '  +7.98%  '

Set-TextValue $ws.Range("D8") '0.613'
$ws.Range("E8").Value = '  +3.81%  '

Set-TextValue $ws.Range("D9") '1.00'
$ws.Range("E9").Value = '  -0.11%  '

Set-TextValue $ws.Range("D10") '0.640'
$ws.Range("E10").Value = '  +8.16%  '

Set-TextValue $ws.Range("D11") '57.29'
$ws.Range("E11").Value = '  +5.11%  '

$ws.Range("E12").Value = '  +15.69%  '

Set-TextValue $ws.Range("D13") '0.0000276'
$ws.Range("E13").Value = '  +10.16%  '

Set-TextValue $ws.Range("D14") '9.50'
$ws.Range("E14").Value = '  +7.53%  '

Set-TextValue $ws.Range("D15") '4.063.60'
$ws.Range("E15").Value = '  +8.35%  '

Set-TextValue $ws.Range("D16") '3.495.91'
$ws.Range("E16").Value = '  +7.97%  '

Set-TextValue $ws.Range("D17") '68.273.82'
$ws.Range("E17").Value = '  +8.27%  '

$ws.Range("E18").Value = '  +6.30%  '

Set-TextValue $ws.Range("D19") '18.36'
$ws.Range("E19").Value = '  +8.05%  '

Set-TextValue $ws.Range("D20") '11.85'
$ws.Range("E20").Value = '  +10.22%  '

$ws.Range("E21").Value = '  +7.72%  '

Set-TextValue $ws.Range("D22") '406.27'
$ws.Range("E22").Value = '  +11.59%  '

Set-TextValue $ws.Range("D23") '11.87'
$ws.Range("E23").Value = '  +11.24%  '

$ws.Range("E24").Value = '  +7.75%  '

Set-TextValue $ws.Range("D25") '84.90'
$ws.Range("E25").Value = '  +7.92%  '

Set-TextValue $ws.Range("D26") '4.22'
$ws.Range("E26").Value = '  +10.82%  '

Set-TextValue $ws.Range("D27") '2.91'
$ws.Range("E27").Value = '  +11.58%  '

Set-TextValue $ws.Range("D28") '6.13'
$ws.Range("E28").Value = '  +1.15%  '

$ws.Range("E29").Value = '  +6.98%  '

Set-TextValue $ws.Range("D30") '8.64'
$ws.Range("E30").Value = '  +7.05%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D31") '30.50'
$ws.Range("E31").Value = '  +8.52%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D32") '687.77'
$ws.Range("E32").Value = '  +9.71%  '

Set-TextValue $ws.Range("D33") '6.88'
$ws.Range("E33").Value = '  +7.52%  '

Set-TextValue $ws.Range("D34") '11.71'
$ws.Range("E34").Value = '  +6.48%  '

$ws.Range("E35").Value = '  +8.08%  '

Set-TextValue $ws.Range("D36") '60.70'
$ws.Range("E36").Value = '  +4.37%  '

Set-TextValue $ws.Range("D37") '39.07'
$ws.Range("E37").Value = '  +9.03%  '

Set-TextValue $ws.Range("D38") '0.0₃0826'
$ws.Range("E38").Value = '  +23.51%  '

$ws.Range("E39").Value = '  +8.10%  '

Set-TextValue $ws.Range("D40") '1.00'
$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D41") '2.80'
$ws.Range("E41").Value = '  +18.12%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D42") '3.37'
$ws.Range("E42").Value = '  +25.34%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D43") '0.134'
$ws.Range("E43").Value = '  +11.91%  '

$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D45") '3.057.46'
$ws.Range("E45").Value = '  +6.70%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range("D46") '2.99'
$ws.Range("E46").Value = '  +14.90%  '

Set-TextValue $ws.Range("D47") '0.0421'
$ws.Range("E47").Value = '  +10.00%  '

$ws.Range("E48").Value = '  +6.40%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D49") '9.21'
$ws.Range("E49").Value = '  +21.23%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D50") '3.23'
$ws.Range("E50").Value = '  +10.98%  '

$ws.Range("E51").Value = '  +6.54%  '
